$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two time-slot labels in column B for rows 20 and 21
$ws.Range("B20").Value = "11:55 - 11:59"
$ws.Range("B21").Value = "12:00 - 12:04"

# Update the trigger JSON values in column C (rows 8 through 21) to the
# corrected trigger list (second entry "04" -> "05")
$ws.Range("C8:C21").Value = '["01", "05", "08", "13", "17", "22", "26", "30"]'

# Update the active selection to match the saved view state (B23)
[void]$ws.Range("B23").Select()
